# Daily attendance processing - 2025-12-31 15:57:39
# Normalizes the "Recorded By" (column G) contributor lists so that
# "System" is reported consistently with the rest of the recorders,
# re-ordering the comma separated list of recorders for each session row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact "before" recorder strings to their normalized "after" strings.
$recorderMap = @{
    "dnasr281@gmail.com, System"            = "System, dnasr281@gmail.com"
    "admin@admin.com, System"               = "System, admin@admin.com"
    "dnasr281@gmail.com, admin@admin.com"   = "admin@admin.com, dnasr281@gmail.com"
    "system, System, backup@backdoor.com"   = "System, backup@backdoor.com, system"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2
    if ($null -ne $value -and $recorderMap.ContainsKey($value)) {
        $cell.Value = $recorderMap[$value]
    }
}
